$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns A, B, E, F, G, H between row 2 and row 3.
$cols = @("A", "B", "E", "F", "G", "H")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")

    $valRow2 = $cellRow2.Value()
    $valRow3 = $cellRow3.Value()

    $cellRow2.Value = $valRow3
    $cellRow3.Value = $valRow2
}
